$wb = $excel.ActiveWorkbook

# Sheet 1
$ws = $wb.Worksheets.Item(1)
$ws.Range("A2").Value = "07-05-2025 21:00"
$ws.Range("B2").Value = "BRAZIL"
$ws.Range("C2").Value = "BRASILEIRO U20 A"
$ws.Range("D2").Value = "Athletico PR U20 - Corinthians U20"
$ws.Range("E2").Value = 73.3
$ws.Range("F2").Value = 1.91
$ws.Range("A3").Value = "08-05-2025 01:00"
$ws.Range("B3").Value = "BRAZIL"
$ws.Range("C3").Value = "MARANHENSE"
$ws.Range("D3").Value = "IAPE - Maranhão"
$ws.Range("E3").Value = 80
$ws.Range("F3").Value = 4.2
$ws.Range("A4").Value = "07-05-2025 18:00"
$ws.Range("B4").Value = "MOROCCO"
$ws.Range("C4").Value = "BOTOLA 2"
$ws.Range("D4").Value = "Olympique Khouribga - Chabab Ben Guerir"
$ws.Range("E4").Value = 73.3
$ws.Range("F4").Value = 2.2
$ws.Range("A5").Value = "08-05-2025 22:00"
$ws.Range("B5").Value = "WORLD"
$ws.Range("C5").Value = "UEFA EUROPA LEAGUE"
$ws.Range("D5").Value = "Bodo/Glimt - Tottenham"
$ws.Range("E5").Value = 90
$ws.Range("F5").Value = 2.9
$ws.Range("A6").Value = "08-05-2025 03:00"
$ws.Range("B6").Value = "WORLD"
$ws.Range("C6").Value = "CONMEBOL LIBERTADORES"
$ws.Range("D6").Value = "Universidad De Chile - Estudiantes L.P."
$ws.Range("E6").Value = 80
$ws.Range("F6").Value = 2.3
$ws.Range("A7").Value = "08-05-2025 05:00"
$ws.Range("B7").Value = "WORLD"
$ws.Range("C7").Value = "CONMEBOL LIBERTADORES"
$ws.Range("D7").Value = "Sporting Cristal - Bolívar"
$ws.Range("E7").Value = 70
$ws.Range("F7").Value = 2.05
$ws.Range("A8").Value = "08-05-2025 03:30"
$ws.Range("B8").Value = "WORLD"
$ws.Range("C8").Value = "CONMEBOL SUDAMERICANA"
$ws.Range("D8").Value = "Mushuc Runa SC - Cruzeiro"
$ws.Range("E8").Value = 70
$ws.Range("F8").Value = 2.2

# Sheet 2
$ws = $wb.Worksheets.Item(2)
$ws.Range("A2").Value = "07-05-2025 20:00"
$ws.Range("B2").Value = "NORWAY"
$ws.Range("C2").Value = "NM CUPEN"
$ws.Range("D2").Value = "Ranheim - Egersund"
$ws.Range("E2").Value = 80
$ws.Range("F2").Value = 2.6

# Sheet 3
$ws = $wb.Worksheets.Item(3)
$ws.Range("A2").Value = "08-05-2025 01:00"
$ws.Range("B2").Value = "BRAZIL"
$ws.Range("C2").Value = "COPA ESPÍRITO SANTO"
$ws.Range("D2").Value = "Real Noroeste - Porto Vitória"
$ws.Range("E2").Value = 83.3
$ws.Range("F2").Value = 1.8
$ws.Range("A3").Value = "07-05-2025 17:00"
$ws.Range("B3").Value = "BULGARIA"
$ws.Range("C3").Value = "FIRST LEAGUE"
$ws.Range("D3").Value = "CSKA 1948 - Hebar 1918"
$ws.Range("E3").Value = 81.1
$ws.Range("F3").Value = 2.02
$ws.Range("A4").Value = "07-05-2025 18:30"
$ws.Range("B4").Value = "CZECH-REPUBLIC"
$ws.Range("C4").Value = "FNL"
$ws.Range("D4").Value = "Opava - Zlin"
$ws.Range("E4").Value = 75
$ws.Range("F4").Value = 1.91
$ws.Range("A5").Value = "07-05-2025 18:00"
$ws.Range("B5").Value = "MOROCCO"
$ws.Range("C5").Value = "BOTOLA 2"
$ws.Range("D5").Value = "Racing De Casablanca - Olympique Dcheïra"
$ws.Range("E5").Value = 86.7
$ws.Range("F5").Value = 1.95
$ws.Range("A6").Value = "07-05-2025 19:00"
$ws.Range("B6").Value = "NORWAY"
$ws.Range("C6").Value = "NM CUPEN"
$ws.Range("D6").Value = "Os - Asane"
$ws.Range("E6").Value = 83.3
$ws.Range("F6").Value = 1.8
$ws.Range("A7").Value = "07-05-2025 19:00"
$ws.Range("B7").Value = "NORWAY"
$ws.Range("C7").Value = "NM CUPEN"
$ws.Range("D7").Value = "Strommen - Mjondalen"
$ws.Range("E7").Value = 78.3
$ws.Range("F7").Value = 1.8
$ws.Range("A8").Value = "07-05-2025 19:00"
$ws.Range("B8").Value = "NORWAY"
$ws.Range("C8").Value = "NM CUPEN"
$ws.Range("D8").Value = "Tromsdalen Uil - Rosenborg"
$ws.Range("E8").Value = 76
$ws.Range("F8").Value = 1.83
$ws.Range("A9").Value = "08-05-2025 03:30"
$ws.Range("B9").Value = "WORLD"
$ws.Range("C9").Value = "CONMEBOL SUDAMERICANA"
$ws.Range("D9").Value = "Atletico Grau - Gremio"
$ws.Range("E9").Value = 88
$ws.Range("F9").Value = 2.1
$ws.Range("A10").Value = "09-05-2025 01:00"
$ws.Range("B10").Value = "WORLD"
$ws.Range("C10").Value = "CONMEBOL SUDAMERICANA"
$ws.Range("D10").Value = "Deportes Iquique - Atletico-MG"
$ws.Range("E10").Value = 88
$ws.Range("F10").Value = 2
$ws.Range("A11").Value = "08-05-2025 20:00"
$ws.Range("B11").Value = "DENMARK"
$ws.Range("C11").Value = "DBU POKALEN"
$ws.Range("D11").Value = "FC Copenhagen - Viborg"
$ws.Range("E11").Value = 76.7
$ws.Range("F11").Value = 1.8

# Sheet 4
$ws = $wb.Worksheets.Item(4)
$ws.Range("A3").Value = "07-05-2025 20:00"
$ws.Range("B3").Value = "NORWAY"
$ws.Range("C3").Value = "NM CUPEN"
$ws.Range("D3").Value = "Ranheim - Egersund"
$ws.Range("E3").Value = 80
$ws.Range("F3").Value = 1.75
$ws.Range("G3").Value = 50
$ws.Range("H3").Value = 2.75
$ws.Range("A4").Value = "07-05-2025 18:30"
$ws.Range("B4").Value = "SLOVENIA"
$ws.Range("C4").Value = "2. SNL"
$ws.Range("D4").Value = "Tabor Sežana - Slovan Ljubljana"
$ws.Range("E4").Value = 73.3
$ws.Range("F4").Value = 1.5
$ws.Range("G4").Value = 60
$ws.Range("H4").Value = 2.4
$ws.Range("A5").Value = "08-05-2025 05:00"
$ws.Range("B5").Value = "WORLD"
$ws.Range("C5").Value = "CONMEBOL LIBERTADORES"
$ws.Range("D5").Value = "Sporting Cristal - Bolívar"
$ws.Range("E5").Value = 86.7
$ws.Range("F5").Value = 1.62
$ws.Range("G5").Value = 60
$ws.Range("H5").Value = 2.38
$ws.Range("A6").Value = "09-05-2025 01:00"
$ws.Range("B6").Value = "WORLD"
$ws.Range("C6").Value = "CONMEBOL SUDAMERICANA"
$ws.Range("D6").Value = "Deportes Iquique - Atletico-MG"
$ws.Range("E6").Value = 80
$ws.Range("F6").Value = 1.95
$ws.Range("G6").Value = 53.3
$ws.Range("H6").Value = 3.4

# Sheet 5
$ws = $wb.Worksheets.Item(5)
$ws.Range("A2").Value = "07-05-2025 03:30"
$ws.Range("B2").Value = "WORLD"
$ws.Range("C2").Value = "CONMEBOL SUDAMERICANA"
$ws.Range("D2").Value = "Boston River - Independiente"
$ws.Range("E2").Value = 60
$ws.Range("F2").Value = 4.75
$ws.Range("G2").Value = 1.85
$ws.Range("A3").Value = "08-05-2025 01:00"
$ws.Range("B3").Value = "WORLD"
$ws.Range("C3").Value = "CONMEBOL SUDAMERICANA"
$ws.Range("D3").Value = "Puerto Cabello - Vasco DA Gama"
$ws.Range("E3").Value = 60
$ws.Range("F3").Value = 4.2
$ws.Range("G3").Value = 1.52
$ws.Range("A4").Value = "07-05-2025 21:00"
$ws.Range("B4").Value = "BRAZIL"
$ws.Range("C4").Value = "BRASILEIRO U20 A"
$ws.Range("D4").Value = "Atlético GO U20 - Juventude U20"
$ws.Range("E4").Value = 50
$ws.Range("F4").Value = 2.3
$ws.Range("G4").Value = 0.15
$ws.Range("A5").Value = "07-05-2025 21:00"
$ws.Range("B5").Value = "BRAZIL"
$ws.Range("C5").Value = "BRASILEIRO U20 A"
$ws.Range("D5").Value = "Athletico PR U20 - Corinthians U20"
$ws.Range("E5").Value = 73.3
$ws.Range("F5").Value = 1.91
$ws.Range("G5").Value = 0.4
$ws.Range("A6").Value = "08-05-2025 01:00"
$ws.Range("B6").Value = "BRAZIL"
$ws.Range("C6").Value = "MARANHENSE"
$ws.Range("D6").Value = "IAPE - Maranhão"
$ws.Range("E6").Value = 80
$ws.Range("F6").Value = 4.2
$ws.Range("G6").Value = 2.36
$ws.Range("A7").Value = "07-05-2025 03:30"
$ws.Range("B7").Value = "BRAZIL"
$ws.Range("C7").Value = "SERIE B"
$ws.Range("D7").Value = "Athletic Club - Vila Nova"
$ws.Range("E7").Value = 60
$ws.Range("F7").Value = 2.3
$ws.Range("G7").Value = 0.38
$ws.Range("A8").Value = "07-05-2025 19:30"
$ws.Range("B8").Value = "BULGARIA"
$ws.Range("C8").Value = "FIRST LEAGUE"
$ws.Range("D8").Value = "Lokomotiv Sofia - Botev Vratsa"
$ws.Range("E8").Value = 60
$ws.Range("F8").Value = 1.7
$ws.Range("G8").Value = 0.02
$ws.Range("A9").Value = "07-05-2025 14:00"
$ws.Range("B9").Value = "CAMBODIA"
$ws.Range("C9").Value = "HUN SEN CUP"
$ws.Range("D9").Value = "Visakha - Phnom Penh Crown"
$ws.Range("E9").Value = 53.3
$ws.Range("F9").Value = 2.15
$ws.Range("G9").Value = 0.15
$ws.Range("A10").Value = "07-05-2025 11:00"
$ws.Range("B10").Value = "CHINA"
$ws.Range("C10").Value = "LEAGUE TWO"
$ws.Range("D10").Value = "Shangyu Pterosaur - Guangxi Hengchen"
$ws.Range("E10").Value = 66.7
$ws.Range("F10").Value = 3.8
$ws.Range("G10").Value = 1.53
$ws.Range("A11").Value = "07-05-2025 19:00"
$ws.Range("B11").Value = "CYPRUS"
$ws.Range("C11").Value = "CUP"
$ws.Range("D11").Value = "AEK Larnaca - Omonia Nicosia"
$ws.Range("E11").Value = 53.3
$ws.Range("F11").Value = 2
$ws.Range("G11").Value = 0.07
$ws.Range("A12").Value = "07-05-2025 18:30"
$ws.Range("B12").Value = "CZECH-REPUBLIC"
$ws.Range("C12").Value = "FNL"
$ws.Range("D12").Value = "Chrudim - Prostějov"
$ws.Range("E12").Value = 60
$ws.Range("F12").Value = 2
$ws.Range("G12").Value = 0.2
$ws.Range("A13").Value = "07-05-2025 18:00"
$ws.Range("B13").Value = "CZECH-REPUBLIC"
$ws.Range("C13").Value = "FNL"
$ws.Range("D13").Value = "Zbrojovka Brno - Viktoria Žižkov"
$ws.Range("E13").Value = 51.7
$ws.Range("F13").Value = 1.73
$ws.Range("G13").Value = -0.11
$ws.Range("A14").Value = "07-05-2025 17:00"
$ws.Range("B14").Value = "EGYPT"
$ws.Range("C14").Value = "PREMIER LEAGUE"
$ws.Range("D14").Value = "El Gouna FC - Future FC"
$ws.Range("E14").Value = 51.3
$ws.Range("F14").Value = 3.35
$ws.Range("G14").Value = 0.72
$ws.Range("A15").Value = "07-05-2025 18:00"
$ws.Range("B15").Value = "ETHIOPIA"
$ws.Range("C15").Value = "PREMIER LEAGUE"
$ws.Range("D15").Value = "Ethiopian Medhin - Bahardar"
$ws.Range("E15").Value = 66.7
$ws.Range("F15").Value = 2.15
$ws.Range("G15").Value = 0.43
$ws.Range("A16").Value = "07-05-2025 18:00"
$ws.Range("B16").Value = "GEORGIA"
$ws.Range("C16").Value = "EROVNULI LIGA"
$ws.Range("D16").Value = "Dila - Torpedo Kutaisi"
$ws.Range("E16").Value = 53.3
$ws.Range("F16").Value = 1.91
$ws.Range("G16").Value = 0.02
$ws.Range("A17").Value = "07-05-2025 18:00"
$ws.Range("B17").Value = "MOROCCO"
$ws.Range("C17").Value = "BOTOLA 2"
$ws.Range("D17").Value = "Olympique Khouribga - Chabab Ben Guerir"
$ws.Range("E17").Value = 73.3
$ws.Range("F17").Value = 2.2
$ws.Range("G17").Value = 0.61
$ws.Range("A18").Value = "07-05-2025 19:00"
$ws.Range("B18").Value = "NORWAY"
$ws.Range("C18").Value = "NM CUPEN"
$ws.Range("D18").Value = "Os - Asane"
$ws.Range("E18").Value = 60
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 6.8
$ws.Range("A19").Value = "07-05-2025 19:00"
$ws.Range("B19").Value = "NORWAY"
$ws.Range("C19").Value = "NM CUPEN"
$ws.Range("D19").Value = "Alta - Kristiansund BK"
$ws.Range("E19").Value = 50
$ws.Range("F19").Value = 5
$ws.Range("G19").Value = 1.5
$ws.Range("A20").Value = "07-05-2025 19:00"
$ws.Range("B20").Value = "NORWAY"
$ws.Range("C20").Value = "NM CUPEN"
$ws.Range("D20").Value = "Rana - Molde"
$ws.Range("E20").Value = 60
$ws.Range("F20").Value = 15
$ws.Range("G20").Value = 8
$ws.Range("A21").Value = "07-05-2025 19:00"
$ws.Range("B21").Value = "NORWAY"
$ws.Range("C21").Value = "NM CUPEN"
$ws.Range("D21").Value = "Pors Grenland - Fredrikstad"
$ws.Range("E21").Value = 50
$ws.Range("F21").Value = 11
$ws.Range("G21").Value = 4.5
$ws.Range("A22").Value = "07-05-2025 18:00"
$ws.Range("B22").Value = "WORLD"
$ws.Range("C22").Value = "AFRICA CUP OF NATIONS U20"
$ws.Range("D22").Value = "Nigeria U20 - Kenya U20"
$ws.Range("E22").Value = 50
$ws.Range("F22").Value = 1.75
$ws.Range("G22").Value = -0.12
$ws.Range("A23").Value = "08-05-2025 22:00"
$ws.Range("B23").Value = "WORLD"
$ws.Range("C23").Value = "UEFA EUROPA LEAGUE"
$ws.Range("D23").Value = "Bodo/Glimt - Tottenham"
$ws.Range("E23").Value = 90
$ws.Range("F23").Value = 2.9
$ws.Range("G23").Value = 1.61
$ws.Range("A24").Value = "08-05-2025 03:00"
$ws.Range("B24").Value = "WORLD"
$ws.Range("C24").Value = "CONMEBOL LIBERTADORES"
$ws.Range("D24").Value = "Universidad De Chile - Estudiantes L.P."
$ws.Range("E24").Value = 80
$ws.Range("F24").Value = 2.3
$ws.Range("G24").Value = 0.84
$ws.Range("A25").Value = "08-05-2025 05:00"
$ws.Range("B25").Value = "WORLD"
$ws.Range("C25").Value = "CONMEBOL LIBERTADORES"
$ws.Range("D25").Value = "Sporting Cristal - Bolívar"
$ws.Range("E25").Value = 70
$ws.Range("F25").Value = 2.05
$ws.Range("G25").Value = 0.43
$ws.Range("A26").Value = "09-05-2025 01:00"
$ws.Range("B26").Value = "WORLD"
$ws.Range("C26").Value = "CONMEBOL LIBERTADORES"
$ws.Range("D26").Value = "Velez Sarsfield - Olimpia"
$ws.Range("E26").Value = 50
$ws.Range("F26").Value = 1.73
$ws.Range("G26").Value = -0.14
$ws.Range("A27").Value = "08-05-2025 03:30"
$ws.Range("B27").Value = "WORLD"
$ws.Range("C27").Value = "CONMEBOL SUDAMERICANA"
$ws.Range("D27").Value = "Mushuc Runa SC - Cruzeiro"
$ws.Range("E27").Value = 70
$ws.Range("F27").Value = 2.2
$ws.Range("G27").Value = 0.54
$ws.Range("A28:G28").EntireRow.Delete()

# Sheet 6
$ws = $wb.Worksheets.Item(6)
$ws.Range("A2").Value = "07-05-2025 03:30"
$ws.Range("B2").Value = "WORLD"
$ws.Range("C2").Value = "CONMEBOL LIBERTADORES"
$ws.Range("D2").Value = "Fortaleza EC - Colo Colo"
$ws.Range("E2").Value = 66.7
$ws.Range("F2").Value = 4.1
$ws.Range("G2").Value = 1.73
$ws.Range("A3").Value = "07-05-2025 10:00"
$ws.Range("B3").Value = "CHINA"
$ws.Range("C3").Value = "LEAGUE TWO"
$ws.Range("D3").Value = "BIT - Nantong Haimen Codion"
$ws.Range("E3").Value = 50
$ws.Range("F3").Value = 1.85
$ws.Range("G3").Value = -0.07
$ws.Range("A4").Value = "07-05-2025 11:00"
$ws.Range("B4").Value = "CZECH-REPUBLIC"
$ws.Range("C4").Value = "1. LIGA U19"
$ws.Range("D4").Value = "Viktoria Plzeň U19 - Baník Ostrava U19"
$ws.Range("E4").Value = 57.7
$ws.Range("F4").Value = 1.92
$ws.Range("G4").Value = 0.11
$ws.Range("A5").Value = "07-05-2025 18:30"
$ws.Range("B5").Value = "CZECH-REPUBLIC"
$ws.Range("C5").Value = "FNL"
$ws.Range("D5").Value = "Opava - Zlin"
$ws.Range("E5").Value = 58.3
$ws.Range("F5").Value = 2.3
$ws.Range("G5").Value = 0.34
$ws.Range("A6").Value = "07-05-2025 18:00"
$ws.Range("B6").Value = "MACEDONIA"
$ws.Range("C6").Value = "FIRST LEAGUE"
$ws.Range("D6").Value = "Akademija Pandev - Shkupi 1927"
$ws.Range("E6").Value = 53.3
$ws.Range("F6").Value = 4.4
$ws.Range("G6").Value = 1.35
$ws.Range("A7").Value = "07-05-2025 20:00"
$ws.Range("B7").Value = "NORWAY"
$ws.Range("C7").Value = "NM CUPEN"
$ws.Range("D7").Value = "Ranheim - Egersund"
$ws.Range("E7").Value = 80
$ws.Range("F7").Value = 2.6
$ws.Range("G7").Value = 1.08
$ws.Range("A8").Value = "07-05-2025 19:00"
$ws.Range("B8").Value = "NORWAY"
$ws.Range("C8").Value = "NM CUPEN"
$ws.Range("D8").Value = "Lyn - Ham-Kam"
$ws.Range("E8").Value = 50
$ws.Range("F8").Value = 2.3
$ws.Range("G8").Value = 0.15
$ws.Range("A9").Value = "07-05-2025 18:00"
$ws.Range("B9").Value = "WORLD"
$ws.Range("C9").Value = "AFRICA CUP OF NATIONS U20"
$ws.Range("D9").Value = "Tunisia U20 - Morocco U20"
$ws.Range("E9").Value = 60
$ws.Range("F9").Value = 2.25
$ws.Range("G9").Value = 0.35
$ws.Range("A10").Value = "08-05-2025 03:30"
$ws.Range("B10").Value = "WORLD"
$ws.Range("C10").Value = "CONMEBOL LIBERTADORES"
$ws.Range("D10").Value = "Cerro Porteno - Palmeiras"
$ws.Range("E10").Value = 66.7
$ws.Range("F10").Value = 1.85
$ws.Range("G10").Value = 0.23
$ws.Range("A11").Value = "09-05-2025 01:00"
$ws.Range("B11").Value = "WORLD"
$ws.Range("C11").Value = "CONMEBOL SUDAMERICANA"
$ws.Range("D11").Value = "Racing Montevideo - Huracan"
$ws.Range("E11").Value = 50
$ws.Range("F11").Value = 2
$ws.Range("G11").Value = 0
$ws.Range("A12").Value = "08-05-2025 18:15"
$ws.Range("B12").Value = "OMAN"
$ws.Range("C12").Value = "PROFESSIONAL LEAGUE"
$ws.Range("D12").Value = "Bahla - Al-Shabab"
$ws.Range("E12").Value = 53.3
$ws.Range("F12").Value = 2.8
$ws.Range("G12").Value = 0.49
$ws.Range("A13:G15").EntireRow.Delete()

# Sheet 7
$ws = $wb.Worksheets.Item(7)
$ws.Range("A2").Value = "07-05-2025 20:00"
$ws.Range("B2").Value = "NORWAY"
$ws.Range("C2").Value = "NM CUPEN"
$ws.Range("D2").Value = "Ranheim - Egersund"
$ws.Range("E2").Value = 80
$ws.Range("F2").Value = 1.75
$ws.Range("G2").Value = 0.4
$ws.Range("A3").Value = "09-05-2025 01:00"
$ws.Range("B3").Value = "WORLD"
$ws.Range("C3").Value = "CONMEBOL SUDAMERICANA"
$ws.Range("D3").Value = "Deportes Iquique - Atletico-MG"
$ws.Range("E3").Value = 80
$ws.Range("F3").Value = 1.95
$ws.Range("G3").Value = 0.56
$ws.Range("A4:G4").EntireRow.Delete()

# Sheet 8
$ws = $wb.Worksheets.Item(8)
$ws.Range("A2").Value = "08-05-2025 01:00"
$ws.Range("B2").Value = "WORLD"
$ws.Range("C2").Value = "CONMEBOL SUDAMERICANA"
$ws.Range("D2").Value = "Puerto Cabello - Vasco DA Gama"
$ws.Range("E2").Value = 68
$ws.Range("F2").Value = 2.1
$ws.Range("G2").Value = 0.43
$ws.Range("A3").Value = "08-05-2025 01:00"
$ws.Range("B3").Value = "BRAZIL"
$ws.Range("C3").Value = "COPA ESPÍRITO SANTO"
$ws.Range("D3").Value = "Real Noroeste - Porto Vitória"
$ws.Range("E3").Value = 83.3
$ws.Range("F3").Value = 1.8
$ws.Range("G3").Value = 0.5
$ws.Range("A4").Value = "07-05-2025 03:30"
$ws.Range("B4").Value = "BRAZIL"
$ws.Range("C4").Value = "SERIE B"
$ws.Range("D4").Value = "Athletic Club - Vila Nova"
$ws.Range("E4").Value = 68
$ws.Range("F4").Value = 2.2
$ws.Range("G4").Value = 0.5
$ws.Range("A5").Value = "07-05-2025 17:00"
$ws.Range("B5").Value = "BULGARIA"
$ws.Range("C5").Value = "FIRST LEAGUE"
$ws.Range("D5").Value = "CSKA 1948 - Hebar 1918"
$ws.Range("E5").Value = 81.1
$ws.Range("F5").Value = 2.02
$ws.Range("G5").Value = 0.64
$ws.Range("A6").Value = "07-05-2025 03:00"
$ws.Range("B6").Value = "COLOMBIA"
$ws.Range("C6").Value = "PRIMERA B"
$ws.Range("D6").Value = "Real Cartagena - Barranquilla"
$ws.Range("E6").Value = 70
$ws.Range("F6").Value = 1.95
$ws.Range("G6").Value = 0.36
$ws.Range("A7").Value = "07-05-2025 18:30"
$ws.Range("B7").Value = "CZECH-REPUBLIC"
$ws.Range("C7").Value = "FNL"
$ws.Range("D7").Value = "Opava - Zlin"
$ws.Range("E7").Value = 75
$ws.Range("F7").Value = 1.91
$ws.Range("G7").Value = 0.43
$ws.Range("A8").Value = "07-05-2025 18:30"
$ws.Range("B8").Value = "CZECH-REPUBLIC"
$ws.Range("C8").Value = "FNL"
$ws.Range("D8").Value = "Varnsdorf - Vysočina Jihlava"
$ws.Range("E8").Value = 66.7
$ws.Range("F8").Value = 1.8
$ws.Range("G8").Value = 0.2
$ws.Range("A9").Value = "07-05-2025 18:00"
$ws.Range("B9").Value = "GEORGIA"
$ws.Range("C9").Value = "EROVNULI LIGA"
$ws.Range("D9").Value = "Dila - Torpedo Kutaisi"
$ws.Range("E9").Value = 66.7
$ws.Range("F9").Value = 1.95
$ws.Range("G9").Value = 0.3
$ws.Range("A10").Value = "07-05-2025 13:00"
$ws.Range("B10").Value = "JAPAN"
$ws.Range("C10").Value = "J1 LEAGUE"
$ws.Range("D10").Value = "Shonan Bellmare - Sanfrecce Hiroshima"
$ws.Range("E10").Value = 66.7
$ws.Range("F10").Value = 1.95
$ws.Range("G10").Value = 0.3
$ws.Range("A11").Value = "07-05-2025 18:00"
$ws.Range("B11").Value = "MACEDONIA"
$ws.Range("C11").Value = "FIRST LEAGUE"
$ws.Range("D11").Value = "Vardar Skopje - Struga"
$ws.Range("E11").Value = 70
$ws.Range("F11").Value = 1.85
$ws.Range("G11").Value = 0.29
$ws.Range("A12").Value = "07-05-2025 18:00"
$ws.Range("B12").Value = "MOROCCO"
$ws.Range("C12").Value = "BOTOLA 2"
$ws.Range("D12").Value = "Yacoub El Mansour - Chabab Atl. Khenifra"
$ws.Range("E12").Value = 72
$ws.Range("F12").Value = 2.25
$ws.Range("G12").Value = 0.62
$ws.Range("A13").Value = "07-05-2025 18:00"
$ws.Range("B13").Value = "MOROCCO"
$ws.Range("C13").Value = "BOTOLA 2"
$ws.Range("D13").Value = "Racing De Casablanca - Olympique Dcheïra"
$ws.Range("E13").Value = 86.7
$ws.Range("F13").Value = 1.95
$ws.Range("G13").Value = 0.69
$ws.Range("A14").Value = "07-05-2025 19:00"
$ws.Range("B14").Value = "NORWAY"
$ws.Range("C14").Value = "NM CUPEN"
$ws.Range("D14").Value = "Os - Asane"
$ws.Range("E14").Value = 83.3
$ws.Range("F14").Value = 1.8
$ws.Range("G14").Value = 0.5
$ws.Range("A15").Value = "07-05-2025 19:00"
$ws.Range("B15").Value = "NORWAY"
$ws.Range("C15").Value = "NM CUPEN"
$ws.Range("D15").Value = "Strommen - Mjondalen"
$ws.Range("E15").Value = 78.3
$ws.Range("F15").Value = 1.8
$ws.Range("G15").Value = 0.41
$ws.Range("A16").Value = "07-05-2025 19:00"
$ws.Range("B16").Value = "NORWAY"
$ws.Range("C16").Value = "NM CUPEN"
$ws.Range("D16").Value = "Tromsdalen Uil - Rosenborg"
$ws.Range("E16").Value = 76
$ws.Range("F16").Value = 1.83
$ws.Range("G16").Value = 0.39
$ws.Range("A17").Value = "08-05-2025 03:30"
$ws.Range("B17").Value = "WORLD"
$ws.Range("C17").Value = "CONMEBOL SUDAMERICANA"
$ws.Range("D17").Value = "Atletico Grau - Gremio"
$ws.Range("E17").Value = 88
$ws.Range("F17").Value = 2.1
$ws.Range("G17").Value = 0.85
$ws.Range("A18").Value = "09-05-2025 01:00"
$ws.Range("B18").Value = "WORLD"
$ws.Range("C18").Value = "CONMEBOL SUDAMERICANA"
$ws.Range("D18").Value = "Deportes Iquique - Atletico-MG"
$ws.Range("E18").Value = 88
$ws.Range("F18").Value = 2
$ws.Range("G18").Value = 0.76
$ws.Range("A19").Value = "08-05-2025 20:00"
$ws.Range("B19").Value = "DENMARK"
$ws.Range("C19").Value = "DBU POKALEN"
$ws.Range("D19").Value = "FC Copenhagen - Viborg"
$ws.Range("E19").Value = 76.7
$ws.Range("F19").Value = 1.8
$ws.Range("G19").Value = 0.38
$ws.Range("A20").Value = "08-05-2025 19:00"
$ws.Range("B20").Value = "NORWAY"
$ws.Range("C20").Value = "NM CUPEN"
$ws.Range("D20").Value = "Tromso - KFUM Oslo"
$ws.Range("E20").Value = 70
$ws.Range("F20").Value = 1.83
$ws.Range("G20").Value = 0.28
